$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GSW vs MIN, Game 2 (2025-05-08) - two rows (away @ MIN perspective, home MIN perspective)
$rows = @(
    @{ Row = 18; A = 16; B = "GSW"; C = "MIN"; D = "away"; E = "2025-05-08"; F = "240:00";
       G = 34; H = 76; I = 0.447; J = 9;  K = 32; L = 0.281; M = 16; N = 25; O = 0.64;
       P = 10; Q = 32; R = 42; S = 23; T = 8;  U = 2; V = 17; W = 16; X = 93;  Y = -24;
       Z = 15; AA = 24; AB = 26; AC = 28; AD = "L" },
    @{ Row = 19; A = 17; B = "MIN"; C = "GSW"; D = "home"; E = "2025-05-08"; F = "240:00";
       G = 44; H = 87; I = 0.506; J = 16; K = 37; L = 0.432; M = 13; N = 20; O = 0.65;
       P = 11; Q = 28; R = 39; S = 33; T = 10; U = 5; V = 12; W = 22; X = 117; Y = 24;
       Z = 29; AA = 27; AB = 29; AC = 32; AD = "W" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A carries the bordered/bold "index" style used by all data rows -
    # copy formatting from the row above so the new cell reuses the same style.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D

    # DATE column: values like "2025-05-08" look like dates, so Excel would
    # normally auto-convert them to date serials. Use a leading apostrophe to
    # force text entry, then reset the style to Normal so no quote-prefix
    # formatting sticks around (matches plain inline/shared-string text cells).
    $ws.Range("E$row").Value = "'" + $r.E
    $ws.Range("E$row").Style = "Normal"

    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
    $ws.Range("W$row").Value = $r.W
    $ws.Range("X$row").Value = $r.X
    $ws.Range("Y$row").Value = $r.Y
    $ws.Range("Z$row").Value = $r.Z
    $ws.Range("AA$row").Value = $r.AA
    $ws.Range("AB$row").Value = $r.AB
    $ws.Range("AC$row").Value = $r.AC
    $ws.Range("AD$row").Value = $r.AD
}
